$d = $word.ActiveDocument

# 1) "4 insightful" -> "4+ insightful"
$d.Content.Find.Execute("4 insightful data nuances", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "4+ insightful data nuances", 2)

# 2) "Your team will present" -> "You will present"
$d.Content.Find.Execute("presentation.  Your team will present", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "presentation.  You will present", 2)

# 3) "include code and PowerPoint slides." -> "include code, bulleted written document and PowerPoint slides with narration."
$d.Content.Find.Execute("submission will include code and PowerPoint slides.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "submission will include code, bulleted written document and PowerPoint slides with narration.", 2)
